$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# --- Update row 12: ac_v value and description ---
# C12: 1 -> 3.28
$ws.Range("C12").Value2 = 3.28
# E12: "Aerodynamic center from root" -> "Aerodynamic center from datum [z]"
$ws.Range("E12").Value2 = "Aerodynamic center from datum [z]"

# --- Insert a new row at 22 for the new "z_h" variable ---
$ws.Rows.Item(22).Insert()

$ws.Range("B22").Value2 = "z_h"
$ws.Range("C22").Value2 = 8.65
$ws.Range("C22").Interior.Color = 65535
$ws.Range("D22").Value2 = "m"
$ws.Range("E22").Value2 = "Measured from datum [z]"

# --- Update the view: scroll position & selection ---
$ws.Activate()
$ws.Range("E18").Select()
$excel.ActiveWindow.ScrollRow = 6
